{"js": "// Replace the date heading and each \"dividend\u00f7divisor=\" expression in the table\n// with the new values, preserving existing run formatting (font, size, etc.).\nconst replacements = [\n  [\"2026-01-03 Saturday\", \"2026-01-04 Sunday\"],\n  [\"839\u00f79=\", \"682\u00f74=\"],\n  [\"963\u00f74=\", \"660\u00f78=\"],\n  [\"907\u00f76=\", \"325\u00f77=\"],\n  [\"177\u00f75=\", \"700\u00f72=\"],\n  [\"302\u00f74=\", \"260\u00f73=\"],\n  [\"297\u00f75=\", \"995\u00f75=\"],\n  [\"668\u00f79=\", \"433\u00f76=\"],\n  [\"250\u00f78=\", \"321\u00f75=\"],\n  [\"422\u00f72=\", \"511\u00f75=\"],\n  [\"564\u00f73=\", \"872\u00f74=\"],\n  [\"927\u00f75=\", \"973\u00f77=\"],\n  [\"294\u00f76=\", \"844\u00f78=\"],\n  [\"190\u00f77=\", \"362\u00f73=\"],\n  [\"557\u00f77=\", \"817\u00f72=\"],\n  [\"509\u00f74=\", \"815\u00f78=\"],\n  [\"334\u00f72=\", \"396\u00f72=\"],\n  [\"642\u00f78=\", \"117\u00f78=\"],\n  [\"345\u00f79=\", \"394\u00f72=\"],\n  [\"562\u00f73=\", \"838\u00f76=\"],\n  [\"520\u00f78=\", \"271\u00f79=\"],\n  [\"813\u00f76=\", \"113\u00f75=\"],\n  [\"119\u00f72=\", \"895\u00f76=\"],\n  [\"694\u00f77=\", \"428\u00f72=\"],\n  [\"198\u00f74=\", \"629\u00f78=\"],\n  [\"168\u00f73=\", \"228\u00f74=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and each \"dividend\u00f7divisor=\" expression in the table\n# with the new values, preserving existing run formatting (font, size, etc.).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"2026-01-03 Saturday\"; new=\"2026-01-04 Sunday\"},\n    @{old=\"839\u00f79=\"; new=\"682\u00f74=\"},\n    @{old=\"963\u00f74=\"; new=\"660\u00f78=\"},\n    @{old=\"907\u00f76=\"; new=\"325\u00f77=\"},\n    @{old=\"177\u00f75=\"; new=\"700\u00f72=\"},\n    @{old=\"302\u00f74=\"; new=\"260\u00f73=\"},\n    @{old=\"297\u00f75=\"; new=\"995\u00f75=\"},\n    @{old=\"668\u00f79=\"; new=\"433\u00f76=\"},\n    @{old=\"250\u00f78=\"; new=\"321\u00f75=\"},\n    @{old=\"422\u00f72=\"; new=\"511\u00f75=\"},\n    @{old=\"564\u00f73=\"; new=\"872\u00f74=\"},\n    @{old=\"927\u00f75=\"; new=\"973\u00f77=\"},\n    @{old=\"294\u00f76=\"; new=\"844\u00f78=\"},\n    @{old=\"190\u00f77=\"; new=\"362\u00f73=\"},\n    @{old=\"557\u00f77=\"; new=\"817\u00f72=\"},\n    @{old=\"509\u00f74=\"; new=\"815\u00f78=\"},\n    @{old=\"334\u00f72=\"; new=\"396\u00f72=\"},\n    @{old=\"642\u00f78=\"; new=\"117\u00f78=\"},\n    @{old=\"345\u00f79=\"; new=\"394\u00f72=\"},\n    @{old=\"562\u00f73=\"; new=\"838\u00f76=\"},\n    @{old=\"520\u00f78=\"; new=\"271\u00f79=\"},\n    @{old=\"813\u00f76=\"; new=\"113\u00f75=\"},\n    @{old=\"119\u00f72=\"; new=\"895\u00f76=\"},\n    @{old=\"694\u00f77=\"; new=\"428\u00f72=\"},\n    @{old=\"198\u00f74=\"; new=\"629\u00f78=\"},\n    @{old=\"168\u00f73=\"; new=\"228\u00f74=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
